# Update the "想去人数" (attendee interest count) figures in column F
# for the 展览 (Exhibition) and 全部类型 (All Types) sheets, which carry
# identical data tables.
$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 2288
    3 = 1742
    4 = 338
    5 = 1096
    6 = 875
    8 = 5855
    9 = 90
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
